{"js": "// 1) Remove the \"_GoBack\" bookmark (w:bookmarkStart/w:bookmarkEnd) that sat\n//    right after \"Underkapitlene skrives slik det er beskrevet i malen.\"\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// 2) Split the run ending in \"...skjerming som opph\u00f8rer etc.. \" into three\n//    runs wrapped with proofErr gramStart/gramEnd markers around \"etc..\",\n//    matching Word's grammar-checker markup:\n//      \"...opph\u00f8rer \" + <proofErr gramStart/> + \"etc..\" + <proofErr gramEnd/> + \" \"\nconst body = context.document.body;\nconst etcResults = body.search(\"etc.. \", { matchCase: true });\netcResults.load(\"items\");\nawait context.sync();\n\nif (etcResults.items.length > 0) {\n  const etcRange = etcResults.items[0];\n  const ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:proofErr w:type=\"gramStart\"/><w:r><w:t>etc..</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t xml:space=\"preserve\"> </w:t></w:r></w:p></w:body></w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n  etcRange.insertOoxml(ooxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3) Add a new empty paragraph at the very end of the body (after the last\n//    paragraph, before the section break).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Remove the hidden \"_GoBack\" bookmark that sat right after\n#    \"Underkapitlene skrives slik det er beskrevet i malen.\"\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2) Split the run ending in \"...skjerming som opph\u00f8rer etc.. \" into three\n#    runs wrapped with proofErr gramStart/gramEnd markers around \"etc..\",\n#    matching Word's grammar-checker markup. We read the paragraph's own\n#    text back (so accented characters are preserved exactly) and only\n#    splice in the proofErr/run boundaries.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Et arkivuttrekk*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $full = $target.Range.Text.TrimEnd([char]13)\n    $marker = \"etc.. \"\n    $idx = $full.IndexOf($marker)\n    if ($idx -ge 0) {\n        $before = $full.Substring(0, $idx)\n\n        $ooxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' + \"`n\" + `\n            '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + \"`n\" + `\n            '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' + \"`n\" + `\n            '<pkg:xmlData>' + \"`n\" + `\n            '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t xml:space=\"preserve\">' + $before + '</w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:t>etc..</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t xml:space=\"preserve\"> </w:t></w:r></w:p></w:body></w:document>' + \"`n\" + `\n            '</pkg:xmlData>' + \"`n\" + `\n            '</pkg:part>' + \"`n\" + `\n            '</pkg:package>'\n\n        $target.Range.InsertXML($ooxml)\n    }\n}\n\n# 3) Add a new empty paragraph at the very end of the body (after the last\n#    paragraph, before the section break).\n$lastParagraph = $d.Paragraphs.Last\n$endRange = $lastParagraph.Range\n$endRange.Collapse(0)\n$endRange.InsertParagraphAfter()\n"}
